$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (theta_se) - standard errors now populated via pickled results instead of (nan)
$ws.Range("B4").Value = "(0.0)"
$ws.Range("C4").Value = "(0.22)"
$ws.Range("D4").Value = "(0.45)"
$ws.Range("E4").Value = "(0.12)"
$ws.Range("F4").Value = "(0.8)"
$ws.Range("G4").Value = "(0.52)"

# Row 6 (lambda_se) - standard errors now populated via pickled results instead of (nan)
$ws.Range("B6").Value = "(0.0)"
$ws.Range("C6").Value = "(0.53)"
$ws.Range("D6").Value = "(0.41)"
$ws.Range("E6").Value = "(0.36)"
$ws.Range("F6").Value = "(0.55)"
$ws.Range("G6").Value = "(0.44)"
